$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hide the detail columns D:O (Volume In/Out, Corte, Contagem ciclica, ILA,
# IRA, Paletes Pendentes, Perdas, HE armazem/transporte, Turn over,
# Absenteismo, Custo MOT) now that the dashboard focuses on later columns.
$ws.Range("D1:O1").EntireColumn.Hidden = $true

# Fill in the newly reported figures for 31/01 (row 193): Largada (Total de
# veiculos), Largada (Veiculos que sairam no horario), Volume and
# Indisponibilidade de Frota.
$ws.Range("P193").Value = 10
$ws.Range("Q193").Value = 10
$ws.Range("R193").Value = 37638.05
$ws.Range("Z193").Value = 0

# Move the active selection to Q186, matching where the user left off.
$ws.Range("Q186").Select()
